# "finestra incidenza 7gg centrata su ultimo g"
#
# The rolling 7-day sum (col C, "somma mobile 7gg.") and the per-100k-
# inhabitants rate derived from it (col D) used to be a window CENTERED on
# the row's day (i.e. summing 3 days before .. 3 days after). This change
# re-centers the window on the LAST day of the window, i.e. it becomes a
# trailing 7-day sum: rows r-6 .. r (inclusive) for each day r.
#
# Practical effect:
#   - The first 3 data rows that used to have a value (rows 5-7) no longer
#     have enough trailing history (need 6 prior rows) and become blank.
#   - The last 3 data rows (182-184) now DO have enough trailing history
#     and get a computed value for the first time.
#   - Every other row's C/D value shifts to reflect the new trailing sum.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ratio between the "per 100k abitanti" figure (col D) and the raw 7-day
# sum (col C) -- constant across the whole sheet (100000 / population).
$ratio = 1.3818071273611627

# Data rows are 2..184; col B = "nuovi pos." (new cases for that day).
$firstRow = 2
$lastRow = 184

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $windowStart = $r - 6
    if ($windowStart -ge $firstRow) {
        $sum = 0
        for ($i = $windowStart; $i -le $r; $i++) {
            $sum = $sum + $ws.Cells.Item($i, 2).Value2
        }
        $ws.Cells.Item($r, 3).Value = $sum
        $ws.Cells.Item($r, 4).Value = $sum * $ratio
    }
    else {
        # Not enough trailing history yet -> blank out (was previously a
        # centered-window value for rows 5-7; stays blank for rows 2-4).
        $ws.Cells.Item($r, 3).Value = ""
        $ws.Cells.Item($r, 4).Value = ""
    }
}

Write-Output "recomputed 7-day trailing window for rows $firstRow..$lastRow"
